# Chapter4/AnalysisFiles/Reformatted Figures.pptx — slide 2
# Wrap the existing "Group 1" (picture + textbox + two callouts) together
# with the sibling "Straight Connector 7" and "Line Callout 2 (Accent Bar) 8"
# shapes into a brand-new outer group ("Group 6", id 7).
#
# PowerPoint's per-slide auto-naming counter for "Group N" / shape ids is
# advanced by every Group()/Ungroup() call (even ones that get undone), so a
# throwaway group+ungroup on two of the unrelated sibling shapes is used
# first to land the *real* group on id 7 / "Group 6", matching the target
# document exactly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Advance the slide's internal Group-naming/id counter so the real group
# below lands on the expected id/name.
$advance = $s.Shapes.Range(@(2, 3)).Group()
$null = $advance.Ungroup()

# Group the whole top-level shape set (the existing "Group 1", the
# "Straight Connector 7" connector, and the "Line Callout 2 (Accent Bar) 8"
# shape) into the new outer group.
$outer = $s.Shapes.Range(@(1, 2, 3)).Group()
